$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rank values ---
# Row 8 (Dart trap / 吹箭): rank 3 -> 2
$ws.Range("B8").Value2 = 2
# Row 12 (Pillar / 石柱): rank 2 -> 1
$ws.Range("B12").Value2 = 1

# --- Add new row 14: 暗门 (Secret door) ---
$ws.Range("A14").Value2 = "暗门"

$ws.Range("B14").Value2 = 1
$ws.Range("B14").NumberFormat = "0_);[Red]\(0\)"

$ws.Range("C14").Value2 = 3
$ws.Range("C14").WrapText = $true

$ws.Range("D14").Value2 = "交锋时：选同一行中1张怪物牌替换房间区的1张怪物牌。"
$ws.Range("D14").WrapText = $true

# --- Add new row 15: 交换机 (Trade machine) ---
$ws.Range("A15").Value2 = "交换机"

$ws.Range("B15").Value2 = 1
$ws.Range("B15").NumberFormat = "0_);[Red]\(0\)"

$ws.Range("C15").Value2 = 3
$ws.Range("C15").WrapText = $true

$ws.Range("E15").Value2 = "Trade machine"

$ws.Range("D15").Value2 = "回合结束时在房间区：选本牌前方1行或后方1行的1张怪物牌，替换房间区或手牌的1张怪物牌。<br>`n回合结束时在手牌区：选手牌的1张怪物牌，替换房间区的1张怪物牌。"
$ws.Range("D15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 28.5

# --- Update selection to match post-edit state ---
$ws.Range("D16").Select()
